$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "EURINFL" is renamed to "EUR REAL" -- update every cell that shares that
# string (the row label in A5 and the column header in E1) so the shared
# string table gets edited in place instead of growing a duplicate entry.
$ws.Range("A5").Value = "EUR REAL"
$ws.Range("E1").Value = "EUR REAL"

# The EUR / EUR REAL correlation (C5, mirrored by formula in E3) is updated.
$ws.Range("C5").Value = 0.95

# A new "EUR INFL" asset is added as both a new column header (F1) and a new
# row label (A6), extending the correlation matrix.
$ws.Range("F1").Value = "EUR INFL"
$ws.Range("A6").Value = "EUR INFL"

# New row 6 holds the EUR INFL correlations against USD / EUR / EURUSD /
# EUR REAL / itself.
$ws.Range("B6").Value = 0.3
$ws.Range("C6").Value = -0.1
$ws.Range("D6").Value = 0.43
$ws.Range("E6").Value = -0.09
$ws.Range("F6").Value = 1

# Column F mirrors row 6 via formulas, the same way column E mirrors row 5.
$ws.Range("F2").Formula = "=B6"
$ws.Range("F3").Formula = "=C6"
$ws.Range("F4").Formula = "=D6"
$ws.Range("F5").Formula = "=E6"

# Selection moved to E10 in the saved file.
[void]$ws.Range("E10").Select()
